# Applies the cryptos.xlsx data refresh described in the commit
# "Updated cryptos list on Tue Aug  1 13:38:45 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.864.74"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.98%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.832.21"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.04%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.01"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6891"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9995"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07700"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.88%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3047"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.34"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07807"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.839.83"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.04%  "

$ws.Range("E13").Value = "  -1.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "90.96"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6809"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.87%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.423"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008302"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.891.69"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.92"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.70%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.075.94"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.17%  "

$ws.Range("E21").Value = "  -2.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.449"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.05%  "

$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1480"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.58"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.785"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.22"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.542"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.97%  "

$ws.Range("E30").Value = "  -2.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.145"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.59%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05098"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7804"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.72%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.851"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.34%  "

$ws.Range("E36").Value = "  -3.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.690"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.58%  "

$ws.Range("E38").Value = "  -1.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.221.81"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9531"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "108.88"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.825"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9990"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.617"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.13%  "

$ws.Range("E46").Value = "  -3.58%  "

$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.977.40"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.47%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5157"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.31%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.10"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -9.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.748"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05904"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.36%  "
